# Add "MTTR Sept." and "Failed Changes Sept." metrics columns next to the
# existing "MTTR"/"Failed Changes" columns (which become the "Oct." variants).
#
# Before: ... G=Lead Time Sept.(Days)  H=MTTR            I=Failed Changes   J=1  K=Lead Time for Changes...  L=6 M=5 N=7
# After : ... G=Lead Time Sept.(Days)  H=MTTR Oct.  I=MTTR Sept.  J=Failed Changes Oct.  K=Failed Changes Sept.  L=1  M=Lead Time for Changes...  N=6 O=5 P=7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns: one right after H (MTTR) and one right after the
# (now shifted) former "Failed Changes" column, so the layout becomes
# H, I(new), J(old I), K(new).
$ws.Columns("I").Insert()
$ws.Columns("K").Insert()

# Rename/populate the header row. Order chosen so new shared-string entries
# are created in the same sequence as the target workbook (MTTR Oct.,
# Failed Changes Oct., MTTR Sept., Failed Changes Sept.).
$ws.Range("H1").Value = "MTTR Oct."
$ws.Range("J1").Value = "Failed Changes Oct."
$ws.Range("I1").Value = "MTTR Sept."
$ws.Range("K1").Value = "Failed Changes Sept."

# New "Sept." data columns start at 0 for every data row.
$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("K5").Value = 0

# Row 3's "MTTR Oct." value changed from 1 to 0.
$ws.Range("H3").Value = 0

# Match the workbook's stored column widths for the new/affected columns
# (ColumnWidth gets 5/6 added internally when persisted, so back that out).
$ws.Columns("H").ColumnWidth = 9.54296875 - 0.8333333333333334
$ws.Columns("I").ColumnWidth = 10.36328125 - 0.8333333333333334
$ws.Columns("J").ColumnWidth = 17.26953125 - 0.8333333333333334
$ws.Columns("K").ColumnWidth = 17.26953125 - 0.8333333333333334

# Update the saved selection/active cell.
$null = $ws.Range("K6").Select()
